$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.99999999926525518
$ws.Cells.Item(2, 1).Value = 0.99288149097986911
$ws.Cells.Item(3, 1).Value = 0.9625255804857924
$ws.Cells.Item(4, 1).Value = 0.95093232376826742
$ws.Cells.Item(5, 1).Value = 0.93969686501075222
$ws.Cells.Item(6, 1).Value = 0.91213382605079274
$ws.Cells.Item(7, 1).Value = 0.90898765334873677
$ws.Cells.Item(8, 1).Value = 0.90493821789300988
$ws.Cells.Item(9, 1).Value = 0.89835876166267603
$ws.Cells.Item(10, 1).Value = 0.89269886124648834
$ws.Cells.Item(11, 1).Value = 0.89193047208957876
$ws.Cells.Item(12, 1).Value = 0.89072443754453545
$ws.Cells.Item(13, 1).Value = 0.88683435134348332
$ws.Cells.Item(14, 1).Value = 0.88591405179399874
$ws.Cells.Item(15, 1).Value = 0.88638954641851297
$ws.Cells.Item(16, 1).Value = 0.883882791533944
$ws.Cells.Item(17, 1).Value = 0.88017455765068453
$ws.Cells.Item(18, 1).Value = 0.87906559994163769
$ws.Cells.Item(19, 1).Value = 0.99165820521952797
$ws.Cells.Item(20, 1).Value = 0.98302518988773613
$ws.Cells.Item(21, 1).Value = 0.98162664797451915
$ws.Cells.Item(22, 1).Value = 0.98036213143401185
$ws.Cells.Item(23, 1).Value = 0.94888902894661986
$ws.Cells.Item(24, 1).Value = 0.93586594919691901
$ws.Cells.Item(25, 1).Value = 0.92940857360521212
$ws.Cells.Item(26, 1).Value = 0.90847455144818423
$ws.Cells.Item(27, 1).Value = 0.90362328266686176
$ws.Cells.Item(28, 1).Value = 0.88213491717779324
$ws.Cells.Item(29, 1).Value = 0.8668575595840462
$ws.Cells.Item(30, 1).Value = 0.86028580840800328
$ws.Cells.Item(31, 1).Value = 0.85263071198664098
$ws.Cells.Item(32, 1).Value = 0.85095125237863645
$ws.Cells.Item(33, 1).Value = 0.8504312095656249
